# Updates cryptos list: refresh Price and Volume(1h) columns, and
# restore natural ranking order for two coin pairs that had been
# listed out of sequence (Cronos/Binance-PegBSC-USD and
# Monero/PolygonEcosystemToken).
#
# Price values are assigned with a leading apostrophe so Excel keeps
# them as plain text (matching the source data, which stores prices
# as text strings, including values with trailing zeros or multiple
# "." group separators that would otherwise be altered/rejected if
# interpreted as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'87.513.20"
$ws.Range("E2").Value = "  -4.07%  "
$ws.Range("D3").Value = "'3.028.44"
$ws.Range("E3").Value = "  -6.10%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'208.48"
$ws.Range("E5").Value = "  -4.69%  "
$ws.Range("D6").Value = "'612.52"
$ws.Range("E6").Value = "  -6.17%  "
$ws.Range("D7").Value = "'0.364"
$ws.Range("E7").Value = "  -9.23%  "
$ws.Range("D8").Value = "'0.780"
$ws.Range("E8").Value = "  +11.25%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'3.026.28"
$ws.Range("E10").Value = "  -6.11%  "
$ws.Range("D11").Value = "'0.585"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("E13").Value = "  -9.84%  "
$ws.Range("D14").Value = "'5.22"
$ws.Range("E14").Value = "  -4.42%  "
$ws.Range("D15").Value = "'87.392.45"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "'3.588.90"
$ws.Range("E16").Value = "  -5.92%  "
$ws.Range("D17").Value = "'31.46"
$ws.Range("E17").Value = "  -6.68%  "
$ws.Range("D18").Value = "'3.040.96"
$ws.Range("E18").Value = "  -5.54%  "
$ws.Range("D19").Value = "'3.24"
$ws.Range("E19").Value = "  -4.37%  "
$ws.Range("D20").Value = "'0.0000196"
$ws.Range("E20").Value = "  -14.39%  "
$ws.Range("D21").Value = "'13.09"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").Value = "'415.17"
$ws.Range("E22").Value = "  -6.32%  "
$ws.Range("D23").Value = "'8.02"
$ws.Range("E23").Value = "  -8.28%  "
$ws.Range("D24").Value = "'4.83"
$ws.Range("E24").Value = "  -5.87%  "
$ws.Range("D25").Value = "'5.37"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'11.60"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").Value = "'81.18"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'3.231.95"
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D32").Value = "'7.91"
$ws.Range("E32").Value = "  -7.48%  "
$ws.Range("D33").Value = "'499.84"
$ws.Range("E33").Value = "  -9.50%  "
$ws.Range("D34").Value = "'3.54"
$ws.Range("E34").Value = "  -14.74%  "
$ws.Range("D35").Value = "'6.59"
$ws.Range("E35").Value = "  -7.44%  "
$ws.Range("D36").Value = "'1.77"
$ws.Range("E36").Value = "  -9.59%  "
$ws.Range("D37").Value = "'1.22"
$ws.Range("E37").Value = "  -7.59%  "
$ws.Range("D38").Value = "'21.82"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").Value = "'22.19"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").Value = "'0.128"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D45").Value = "'1.77"
$ws.Range("E45").Value = "  -9.98%  "
$ws.Range("D46").Value = "'0.129"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").Value = "'43.21"
$ws.Range("E47").Value = "  -5.06%  "
$ws.Range("D48").Value = "'0.0675"
$ws.Range("E48").Value = "  +10.46%  "
$ws.Range("D49").Value = "'155.62"
$ws.Range("E49").Value = "  -11.17%  "
$ws.Range("D50").Value = "'1.17"
$ws.Range("E50").Value = "  -6.58%  "
$ws.Range("D51").Value = "'0.690"
$ws.Range("E51").Value = "  -10.69%  "

# Swap rows 30 and 31 (Cronos and Binance-PegBSC-USD swap order)
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").Value = "'0.178"
$ws.Range("E30").Value = "  +9.91%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.08"
$ws.Range("E31").Value = "  +8.83%  "

# Swap rows 43 and 44 (Monero and PolygonEcosystemToken swap order)
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'147.59"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.355"
$ws.Range("E44").Value = "  -6.29%  "
